$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NPC names (row 5 and 6, column B)
$ws.Range("B5").Value = "Tanuki"
$ws.Range("B6").Value = "Frog"

# Row 7: Tanuki's first dialogue line, with new sprite column (D)
$ws.Range("B7").Value = "Tanuki"
$ws.Range("C7").Value = "It'sa me! Tanuki man!"
$ws.Range("D7").Value = "tanuki_mario"

# Row 8: Frog's dialogue line, with new sprite column (D)
$ws.Range("B8").Value = "Frog"
$ws.Range("C8").Value = "Ribbit Ribbit Ribbit Ribbit. Ribbit Ribbit. Riiiiibbbbiiiiiiiiit."
$ws.Range("D8").Value = "frog_mario"

# Row 9: Tanuki's final line, with new sprite column (D)
$ws.Range("B9").Value = "Tanuki"
$ws.Range("C9").Value = "Well, goodbye."
$ws.Range("D9").Value = "tanuki_mario"

# Update the active selection to reflect the last edited cell
$ws.Range("D7").Select()
